$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Cells.Item(1, 5).Value = 'additional_feilds'

# Per-row updates: reformat column D address, add column E unit/owner info
$ws.Cells.Item(2, 4).Value = '36801 Kevin Harbors Apt. 656, Port Laurenport, MD 59547'
$ws.Cells.Item(2, 5).Value = 'Unit ID : asdf123    Owner : Nice   '

$ws.Cells.Item(3, 4).Value = '50393 Jeremiah Ports, Fergusonview, MO 91646'
$ws.Cells.Item(3, 5).Value = 'Unit ID : asdf124    Owner : Nice   '

$ws.Cells.Item(4, 4).Value = '8638 Wilson Overpass, Sheilaborough, IA 91232'
$ws.Cells.Item(4, 5).Value = 'Unit ID : asdf125    Owner : Nice   '

$ws.Cells.Item(5, 4).Value = '627 Jeffrey Valleys, North Tracyfurt, OH 95791'
$ws.Cells.Item(5, 5).Value = 'Unit ID : asdf126    Owner : Nice   '

$ws.Cells.Item(6, 4).Value = '024 Velazquez Lakes, North Samantha, KS 49562'
$ws.Cells.Item(6, 5).Value = 'Unit ID : asdf127    Owner : Nice   '

$ws.Cells.Item(7, 4).Value = '2107 Mathews Mews Apt. 025, Reidmouth, KY 96338'
$ws.Cells.Item(7, 5).Value = 'Unit ID : asdf128    Owner : Nice   '

$ws.Cells.Item(8, 4).Value = '3573 Brown Gardens, Biancahaven, VT 98687'
$ws.Cells.Item(8, 5).Value = 'Unit ID : asdf129    Owner : Nice   '

$ws.Cells.Item(9, 4).Value = '3640 Flores Garden, Lake Yolanda, GA 39886'
$ws.Cells.Item(9, 5).Value = 'Unit ID : asdf130    Owner : Nice   '

$ws.Cells.Item(10, 4).Value = '180 Burke Circle, Hobbschester, ND 71304'
$ws.Cells.Item(10, 5).Value = 'Unit ID : asdf132    Owner : Nice   '

$ws.Cells.Item(11, 4).Value = '81384 Richard View Apt. 816, Emilychester, OR 92406'
$ws.Cells.Item(11, 5).Value = 'Unit ID : asdf133    Owner : Nice   '

$ws.Cells.Item(12, 4).Value = '042 Jamie Hollow, Lisaside, MI 90177'
$ws.Cells.Item(12, 5).Value = 'Unit ID : asdf135    Owner : Not Nice    '

$ws.Cells.Item(13, 4).Value = '36100 Horne Curve Apt. 345, North Lauraton, NH 51356'
$ws.Cells.Item(13, 5).Value = 'Unit ID : asdf136    Owner : Not Nice    '

$ws.Cells.Item(14, 4).Value = '20695 Vaughan Mountain, East Juanfurt, NH 99845'
$ws.Cells.Item(14, 5).Value = 'Unit ID : asdf137    Owner : Not Nice    '

$ws.Cells.Item(15, 4).Value = '6416 Hughes Forks Suite 776, Taylorshire, IA 47964'
$ws.Cells.Item(15, 5).Value = 'Unit ID : asdf138    Owner : Not Nice    '

$ws.Cells.Item(16, 4).Value = '11176 Sierra Greens Suite 643, Deleonfurt, ID 31790'
$ws.Cells.Item(16, 5).Value = 'Unit ID : asdf139    Owner : Not Nice    '

$ws.Cells.Item(17, 4).Value = '988 Vasquez Burgs, Lake Donaldbury, PA 29774'
$ws.Cells.Item(17, 5).Value = 'Unit ID : asdf140    Owner : Not Nice    '

$ws.Cells.Item(18, 4).Value = '91896 Fleming Track, South Stephaniechester, FL 51489'
$ws.Cells.Item(18, 5).Value = 'Unit ID : asdf141    Owner : Not Nice    '

$ws.Cells.Item(19, 4).Value = '02638 Morrison Meadow, Greenbury, MN 54161'
$ws.Cells.Item(19, 5).Value = 'Unit ID : asdf142    Owner : Not Nice    '
